{"js": "// Update the date heading and the 25 division-problem cells in the table.\nconst body = context.document.body;\n\n// --- 1) Date paragraph (first paragraph of the body, before the table) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-04-05 Saturday\", \"Replace\");\n\n// --- 2) Division problems in the table ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Each block of problems lives in rows 0, 4, 8, 12, 16 (5 columns each);\n// the rows in between hold blank answer cells. New values below follow\n// the exact reading order of the source table (left-to-right, top-to-\n// bottom across the five problem rows).\nconst newValues = [\n  [\"94\u00f76=\", \"49\u00f77=\", \"83\u00f77=\", \"10\u00f79=\", \"68\u00f79=\"],\n  [\"78\u00f72=\", \"67\u00f72=\", \"28\u00f73=\", \"15\u00f78=\", \"22\u00f73=\"],\n  [\"27\u00f74=\", \"88\u00f73=\", \"74\u00f72=\", \"75\u00f78=\", \"27\u00f79=\"],\n  [\"57\u00f79=\", \"69\u00f75=\", \"56\u00f77=\", \"75\u00f72=\", \"15\u00f76=\"],\n  [\"16\u00f72=\", \"55\u00f74=\", \"39\u00f75=\", \"72\u00f75=\", \"29\u00f75=\"],\n];\nconst problemRows = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < problemRows.length; i++) {\n  const rowIndex = problemRows[i];\n  for (let col = 0; col < newValues[i].length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-problem cells in the table.\n$d = $word.ActiveDocument\n\n# --- 1) Date paragraph (first paragraph of the document, before the table) ---\n$d.Paragraphs(1).Range.Text = \"2025-04-05 Saturday\"\n\n# --- 2) Division problems in the table ---\n$t = $d.Tables.Item(1)\n\n# Each block of problems lives in (1-indexed) rows 1, 5, 9, 13, 17 (5 columns\n# each); the rows in between hold blank answer cells. New values below\n# follow the exact reading order of the source table (left-to-right,\n# top-to-bottom across the five problem rows).\n$newValues = @(\n    @(\"94\u00f76=\", \"49\u00f77=\", \"83\u00f77=\", \"10\u00f79=\", \"68\u00f79=\"),\n    @(\"78\u00f72=\", \"67\u00f72=\", \"28\u00f73=\", \"15\u00f78=\", \"22\u00f73=\"),\n    @(\"27\u00f74=\", \"88\u00f73=\", \"74\u00f72=\", \"75\u00f78=\", \"27\u00f79=\"),\n    @(\"57\u00f79=\", \"69\u00f75=\", \"56\u00f77=\", \"75\u00f72=\", \"15\u00f76=\"),\n    @(\"16\u00f72=\", \"55\u00f74=\", \"39\u00f75=\", \"72\u00f75=\", \"29\u00f75=\")\n)\n$problemRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $problemRows.Length; $i++) {\n    $rowIndex = $problemRows[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $newValues[$i][$col - 1]\n    }\n}\n"}
